$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same
$ws.Range("A1").Value = "Variable/Konfiguration"
$ws.Range("B1").Value = "Config1"
$ws.Range("C1").Value = "Config2"

# Row 2: replace old generic Variable1/Wert1/Value1 with the new base-state row
$ws.Range("A2").Value = "BaseStateChromium"
$ws.Range("B2").Value = "Started"
$ws.Range("C2").Value = "Not Started"

# Row 3: replace old generic Variable2/Wert2/Value2 with the new base-state row
$ws.Range("A3").Value = "BaseStatePixel9Pro_API35"
$ws.Range("B3").Value = "Not Started"
$ws.Range("C3").Value = "Started"

# Row 4: new row added for AUT configuration
$ws.Range("A4").Value = "AUT"
$ws.Range("B4").Value = "Chromium"
$ws.Range("C4").Value = "Pixel9Pro_API35"

# Update the selection to match the new active selection in the sheet
$ws.Range("A2:XFD4").Select()
